# Apply crypto price/volume updates per commit "Updated cryptos list on Thu Aug 15 07:08:52 UTC 2024 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '58.074.48'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -4.06%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.618.54'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -3.02%  '

$ws.Range("E4").Value = '  -0.16%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '515.72'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.92%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '141.95'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.17%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.999'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.28%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.565'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -2.09%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '6.68'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.29%  '

$ws.Range("E10").Value = '  -3.34%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.335'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.41%  '

$ws.Range("E12").Value = '  +1.30%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '3.076.27'
$ws.Range("D13").Style = "Normal"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '58.063.44'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -4.13%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '20.62'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -3.08%  '

$ws.Range("B16").Value = 'WrappedEther'
$ws.Range("C16").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.631.63'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -3.70%  '

$ws.Range("B17").Value = 'ShibaInu'
$ws.Range("C17").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000135'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -2.00%  '

$ws.Range("B18").Value = 'Polkadot'
$ws.Range("C18").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '4.38'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -2.81%  '

$ws.Range("B19").Value = 'BitcoinCash'
$ws.Range("C19").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '334.41'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -3.30%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '10.30'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -3.15%  '

$ws.Range("E21").Value = '  -3.26%  '

$ws.Range("E22").Value = '  +0.05%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '63.68'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.21%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.421'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.05%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.999'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.66%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.02'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -3.72%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.0₃0779'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -4.86%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '6.55'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -3.87%  '

$ws.Range("E30").Value = '  +0.12%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.58'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.22%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '151.35'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.49%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '18.64'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -2.29%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.06'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -4.37%  '

$ws.Range("B35").Value = 'ImmutableX'
$ws.Range("C35").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.16'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -5.64%  '

$ws.Range("B36").Value = 'SuiNetwork'
$ws.Range("C36").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.891'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -5.18%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '36.58'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.37%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.840'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -3.46%  '

$ws.Range("E39").Value = '  -6.09%  '

$ws.Range("E41").Value = '  +0.37%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.596'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.18%  '

$ws.Range("E43").Value = '  -2.67%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '267.27'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -5.36%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '10.60'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.25%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '19.06'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -5.44%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0530'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.78%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.026.90'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -5.28%  '

$ws.Range("E49").Value = '  -2.70%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '4.61'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -3.94%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '18.19'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -4.87%  '
